$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "[-, -, 'MEC-3B-Tec. Soldagem', -]"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "['MCT-3A-Tecnologia da soldagem', -, -, -]"

# Row 4
$ws.Range("B4").Value = "[-, -, 'MEC-3B-Tec. Soldagem', -]"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "['MCT-3A-Tecnologia da soldagem', -, -, -]"

# Row 6
$ws.Range("B6").Value = "[-, -, 'MEC-3B-Tec. Soldagem', -]"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "['MCT-3A-Tecnologia da soldagem', -, -, -]"

# Row 7
$ws.Range("B7").Value = "[-, -, 'MEC-3B-Tec. Soldagem', -]"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "['MCT-3A-Tecnologia da soldagem', -, -, -]"

# Row 18
$ws.Range("C18").Value = "ELM-1NA-Gestão Integrada"

# Row 19
$ws.Range("F19").Value = "-"

# Row 20
$ws.Range("C20").Value = "MEC-2NA-Gest. Int."
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "MEC-2NB-Gestão integrada"

# Row 21
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "ELM-1NA-Gestão Integrada"
$ws.Range("F21").Value = "MEC-2NB-Gestão integrada"
